$d = $word.ActiveDocument

# Locate the two paragraphs that hold the "Prosecutor's Office..." /
# "{% if community_control... %}" text block by searching for a unique
# anchor string, so this does not depend on a hard-coded paragraph index.
$anchor = $d.Content
$null = $anchor.Find.Execute("Prosecutor" + [char]0x2019 + "s Office:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Expand(4) | Out-Null   # wdParagraph -> expand to the whole 1st paragraph

$nextPara = $d.Range($anchor.End, $anchor.End)
$nextPara.Expand(4) | Out-Null  # the following paragraph (the jinja conditions one)

$fullRng = $d.Range($anchor.Start, $nextPara.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:tabs><w:tab w:val="left" w:pos="-1080"/><w:tab w:val="left" w:pos="-720"/><w:tab w:val="left" w:pos="0"/><w:tab w:val="left" w:pos="720"/><w:tab w:val="left" w:pos="1440"/><w:tab w:val="left" w:pos="2160"/><w:tab w:val="left" w:pos="2880"/><w:tab w:val="left" w:pos="3600"/><w:tab w:val="left" w:pos="4320"/><w:tab w:val="left" w:pos="4680"/></w:tabs><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Prosecutor&#x2019;s Office: PS    OM     EM; Defendant&#x2019;s Attorney: PS     OM     EM; {{ defendant.first_name }} {{ defendant.last_name}}: PS     OM     EM;</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Community Control: PS    EM;</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>County Jail: PS   EM;</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$fullRng.InsertXML($xml)
